$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rng, $val) {
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '40.980.78'
Set-TextValue $ws.Range('E2') '  -1.64%  '
Set-TextValue $ws.Range('D3') '2.414.46'
Set-TextValue $ws.Range('E3') '  -2.41%  '
Set-TextValue $ws.Range('D4') '0.999'
Set-TextValue $ws.Range('E4') '  -0.06%  '
Set-TextValue $ws.Range('D5') '315.69'
Set-TextValue $ws.Range('E5') '  -0.80%  '
Set-TextValue $ws.Range('D6') '88.58'
Set-TextValue $ws.Range('E6') '  -4.50%  '
Set-TextValue $ws.Range('D7') '0.537'
Set-TextValue $ws.Range('E7') '  -2.77%  '
Set-TextValue $ws.Range('E8') '  +0.08%  '
Set-TextValue $ws.Range('D9') '0.495'
Set-TextValue $ws.Range('E9') '  -3.95%  '
Set-TextValue $ws.Range('D10') '0.0834'
Set-TextValue $ws.Range('E10') '  -2.23%  '
Set-TextValue $ws.Range('D11') '31.44'
Set-TextValue $ws.Range('E11') '  -5.00%  '
Set-TextValue $ws.Range('E12') '  -1.68%  '
Set-TextValue $ws.Range('D13') '2.784.89'
Set-TextValue $ws.Range('E13') '  -2.46%  '
Set-TextValue $ws.Range('D14') '6.79'
Set-TextValue $ws.Range('E14') '  -1.52%  '
Set-TextValue $ws.Range('D15') '15.64'
Set-TextValue $ws.Range('E15') '  -0.92%  '
Set-TextValue $ws.Range('D16') '2.435.16'
Set-TextValue $ws.Range('E16') '  -2.97%  '
Set-TextValue $ws.Range('D17') '0.770'
Set-TextValue $ws.Range('E17') '  -2.42%  '
Set-TextValue $ws.Range('D18') '40.845.59'
Set-TextValue $ws.Range('E18') '  -1.89%  '
Set-TextValue $ws.Range('D19') '0.0₃0920'
Set-TextValue $ws.Range('E19') '  -3.34%  '
Set-TextValue $ws.Range('D20') '6.23'
Set-TextValue $ws.Range('E20') '  -3.78%  '
Set-TextValue $ws.Range('D21') '70.79'
Set-TextValue $ws.Range('E21') '  -0.67%  '
Set-TextValue $ws.Range('D22') '10.86'
Set-TextValue $ws.Range('E22') '  -3.99%  '
Set-TextValue $ws.Range('D23') '233.48'
Set-TextValue $ws.Range('E23') '  -2.58%  '
Set-TextValue $ws.Range('D24') '2.67'
Set-TextValue $ws.Range('E24') '  -2.90%  '
Set-TextValue $ws.Range('E25') '  +0.21%  '
Set-TextValue $ws.Range('E26') '  -4.50%  '
Set-TextValue $ws.Range('D27') '23.98'
Set-TextValue $ws.Range('D28') '2.23'
Set-TextValue $ws.Range('E28') '  -2.58%  '
Set-TextValue $ws.Range('D29') '9.53'
Set-TextValue $ws.Range('E29') '  -3.23%  '
Set-TextValue $ws.Range('D30') '34.14'
Set-TextValue $ws.Range('E30') '  -5.59%  '
Set-TextValue $ws.Range('D31') '157.54'
Set-TextValue $ws.Range('E31') '  -1.06%  '
Set-TextValue $ws.Range('B32') 'Filecoin'
Set-TextValue $ws.Range('C32') 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Range('D32') '5.27'
Set-TextValue $ws.Range('E32') '  -4.53%  '
Set-TextValue $ws.Range('B33') 'FirstDigitalUSD'
Set-TextValue $ws.Range('C33') 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue $ws.Range('D33') '1.00'
Set-TextValue $ws.Range('E33') '  -0.03%  '
Set-TextValue $ws.Range('D34') '0.0737'
Set-TextValue $ws.Range('E34') '  -4.01%  '
Set-TextValue $ws.Range('D35') '2.46'
Set-TextValue $ws.Range('E35') '  -4.81%  '
Set-TextValue $ws.Range('D36') '2.88'
Set-TextValue $ws.Range('E36') '  -1.47%  '
Set-TextValue $ws.Range('D37') '16.36'
Set-TextValue $ws.Range('E37') '  -5.41%  '
Set-TextValue $ws.Range('D38') '0.114'
Set-TextValue $ws.Range('E38') '  -1.40%  '
Set-TextValue $ws.Range('B39') 'ARBITRUM'
Set-TextValue $ws.Range('C39') 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws.Range('D39') '1.76'
Set-TextValue $ws.Range('E39') '  -6.02%  '
Set-TextValue $ws.Range('B40') 'Kaspa'
Set-TextValue $ws.Range('C40') 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws.Range('D40') '0.100'
Set-TextValue $ws.Range('E40') '  -3.55%  '
Set-TextValue $ws.Range('D41') '3.87'
Set-TextValue $ws.Range('E41') '  -3.33%  '
Set-TextValue $ws.Range('D42') '2.30'
Set-TextValue $ws.Range('E42') '  -6.67%  '
Set-TextValue $ws.Range('D43') '1.984.29'
Set-TextValue $ws.Range('E43') '  -0.44%  '
Set-TextValue $ws.Range('D44') '18.26'
Set-TextValue $ws.Range('E44') '  -3.89%  '
Set-TextValue $ws.Range('D45') '0.0273'
Set-TextValue $ws.Range('E45') '  -4.28%  '
Set-TextValue $ws.Range('E46') '  -4.72%  '
Set-TextValue $ws.Range('D47') '9.37'
Set-TextValue $ws.Range('E47') '  +0.48%  '
Set-TextValue $ws.Range('D48') '2.645.26'
Set-TextValue $ws.Range('E48') '  -2.47%  '
Set-TextValue $ws.Range('D49') '94.02'
Set-TextValue $ws.Range('E49') '  -3.35%  '
Set-TextValue $ws.Range('D50') '73.22'
Set-TextValue $ws.Range('E50') '  -1.24%  '
Set-TextValue $ws.Range('D51') '51.22'
Set-TextValue $ws.Range('E51') '  -1.92%  '
